$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 244, shifting existing rows 244-281 down to 245-282.
$ws.Rows(244).Insert()

# Populate the newly inserted row 244 with the new weekly record.
$ws.Range("A244").Value = 11
$ws.Range("B244").Value = "Vega Monumental Concepción"
$ws.Range("C244").Value = "Bíobío"
$ws.Range("D244").Value = 44637
$ws.Range("E244").Value = 8
$ws.Range("F244").Value = 100114014
$ws.Range("G244").Value = "Betarraga"
$ws.Range("H244").Value = "Sin especificar"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 200
$ws.Range("K244").Value = 600
$ws.Range("L244").Value = 650
$ws.Range("M244").Value = 625
$ws.Range("N244").Value = "`$/paquete 5 unidades"
$ws.Range("O244").Value = "Región Metropolitana"
$ws.Range("P244").Value = 125
$ws.Range("Q244").Value = 5
$ws.Range("R244").Value = "Hortaliza"

# Match the style used by the rest of column D (date number format).
$ws.Range("D244").NumberFormat = $ws.Range("D245").NumberFormat
